# Update "想去人数" (Column F) values on sheets 展览, 演出 and 全部类型
# to match freshly generated data (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
  2  = 0
  4  = 0
  5  = 4678
  8  = 0
  13 = 1062
  15 = 254
  18 = 0
  20 = 3695
  21 = 6035
  27 = 3407
  28 = 379
  29 = 0
  31 = 0
  33 = 0
  36 = 0
  38 = 0
  40 = 29
  43 = 466
  45 = 0
  46 = 0
}
foreach ($row in $sheet1Updates.Keys) {
  $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

# --- Sheet "演出" (performances) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 0

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
  4  = 0
  5  = 4678
  9  = 101
  10 = 0
  12 = 0
  13 = 0
  14 = 1062
  16 = 254
  18 = 0
  19 = 130
  21 = 0
  22 = 6035
  23 = 0
  25 = 0
  26 = 0
  27 = 0
  28 = 0
  29 = 379
  33 = 0
  38 = 0
  39 = 1539
  40 = 0
  41 = 0
  42 = 0
  43 = 0
  44 = 467
  45 = 0
  46 = 0
}
foreach ($row in $sheet4Updates.Keys) {
  $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
